# Add check for whether a VLAN is a DHCP vlan:
#  - row 2 (10.9.106.30 / floor7_sw_lab): type changes from "access" to "backbone"
#  - row 7 (10.9.107.254 / PT-SW-DIST-Shaked): group changes from "junos" to "none"
#  - row 8 (10.9.106.44 / PTSWCORE2): group changes from "cisco" to "none"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E2").Value = "backbone"
$ws1.Range("B7").Value = "none"
$ws1.Range("B8").Value = "none"

# Make Sheet1 active and move the selected cell to D13
$ws1.Activate()
$ws1.Range("D13").Select()

$wb.Save()
